$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal text, regardless of whether it looks numeric,
# without leaving a lasting NumberFormat change on the cell.
function Set-TextValue($cellRange, $text) {
    $cellRange.NumberFormat = "@"
    $cellRange.Value2 = $text
    $cellRange.ClearFormats()
}

# Row 2
Set-TextValue $ws.Range("D2") '66.423.97'
$ws.Range("E2").Value2 = '  -0.21%  '

# Row 3
Set-TextValue $ws.Range("D3") '3.534.77'
$ws.Range("E3").Value2 = '  -1.82%  '

# Row 5
Set-TextValue $ws.Range("D5") '606.91'
$ws.Range("E5").Value2 = '  -0.43%  '

# Row 6
Set-TextValue $ws.Range("D6") '143.28'
$ws.Range("E6").Value2 = '  -3.62%  '

# Row 7
Set-TextValue $ws.Range("D7") '3.530.29'
$ws.Range("E7").Value2 = '  -1.91%  '

# Row 8
$ws.Range("E8").Value2 = '  +0.00%  '

# Row 9
Set-TextValue $ws.Range("D9") '0.517'
$ws.Range("E9").Value2 = '  +5.55%  '

# Row 10
Set-TextValue $ws.Range("D10") '0.132'
$ws.Range("E10").Value2 = '  -3.40%  '

# Row 11
Set-TextValue $ws.Range("D11") '7.70'
$ws.Range("E11").Value2 = '  -4.75%  '

# Row 12
Set-TextValue $ws.Range("D12") '0.410'
$ws.Range("E12").Value2 = '  -1.53%  '

# Row 13
Set-TextValue $ws.Range("D13") '4.140.55'
$ws.Range("E13").Value2 = '  -1.70%  '

# Row 14
Set-TextValue $ws.Range("D14") '0.0000195'
$ws.Range("E14").Value2 = '  -6.55%  '

# Row 15
Set-TextValue $ws.Range("D15") '28.84'
$ws.Range("E15").Value2 = '  -3.34%  '

# Row 16
Set-TextValue $ws.Range("D16") '3.547.12'
$ws.Range("E16").Value2 = '  -1.50%  '

# Row 17
$ws.Range("E17").Value2 = '  +0.73%  '

# Row 18
Set-TextValue $ws.Range("D18") '66.479.65'
$ws.Range("E18").Value2 = '  -0.30%  '

# Row 19
Set-TextValue $ws.Range("D19") '10.86'
$ws.Range("E19").Value2 = '  -5.50%  '

# Row 20
Set-TextValue $ws.Range("D20") '6.20'
$ws.Range("E20").Value2 = '  -2.84%  '

# Row 21
Set-TextValue $ws.Range("D21") '14.62'
$ws.Range("E21").Value2 = '  -3.24%  '

# Row 22
Set-TextValue $ws.Range("D22") '426.10'
$ws.Range("E22").Value2 = '  -0.46%  '

# Row 23
Set-TextValue $ws.Range("D23") '0.594'
$ws.Range("E23").Value2 = '  -3.99%  '

# Row 24
Set-TextValue $ws.Range("D24") '77.29'
$ws.Range("E24").Value2 = '  -2.05%  '

# Row 25
Set-TextValue $ws.Range("D25") '3.684.02'
$ws.Range("E25").Value2 = '  -1.64%  '

# Row 26
$ws.Range("E26").Value2 = '  -0.11%  '

# Row 27
Set-TextValue $ws.Range("D27") '0.0000115'
$ws.Range("E27").Value2 = '  -6.39%  '

# Row 28
Set-TextValue $ws.Range("D28") '8.00'
$ws.Range("E28").Value2 = '  -4.01%  '

# Row 29
Set-TextValue $ws.Range("D29") '2.47'
$ws.Range("E29").Value2 = '  -2.14%  '

# Row 30
Set-TextValue $ws.Range("D30") '8.98'
$ws.Range("E30").Value2 = '  -5.81%  '

# Row 31
$ws.Range("E31").Value2 = '  +0.19%  '

# Row 32
$ws.Range("B32").Value2 = 'RenzoRestakedETH'
$ws.Range("C32").Value2 = 'https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth'
Set-TextValue $ws.Range("D32") '3.549.30'
$ws.Range("E32").Value2 = '  -1.37%  '

# Row 33
$ws.Range("B33").Value2 = 'Kaspa'
$ws.Range("C33").Value2 = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
Set-TextValue $ws.Range("D33") '0.156'
$ws.Range("E33").Value2 = '  -0.82%  '

# Row 34
Set-TextValue $ws.Range("D34") '24.39'
$ws.Range("E34").Value2 = '  -4.24%  '

# Row 35
$ws.Range("E35").Value2 = '  +0.02%  '

# Row 36
Set-TextValue $ws.Range("D36") '1.35'
$ws.Range("E36").Value2 = '  -8.45%  '

# Row 37
Set-TextValue $ws.Range("D37") '7.57'
$ws.Range("E37").Value2 = '  -3.78%  '

# Row 38
Set-TextValue $ws.Range("D38") '1.63'
$ws.Range("E38").Value2 = '  -4.31%  '

# Row 39
Set-TextValue $ws.Range("D39") '176.03'
$ws.Range("E39").Value2 = '  -0.35%  '

# Row 40
Set-TextValue $ws.Range("D40") '5.25'
$ws.Range("E40").Value2 = '  -7.09%  '

# Row 41
Set-TextValue $ws.Range("D41") '0.0823'
$ws.Range("E41").Value2 = '  -4.00%  '

# Row 42
Set-TextValue $ws.Range("D42") '5.00'
$ws.Range("E42").Value2 = '  -4.79%  '

# Row 43
Set-TextValue $ws.Range("D43") '0.858'
$ws.Range("E43").Value2 = '  -4.56%  '

# Row 44
Set-TextValue $ws.Range("D44") '45.28'
$ws.Range("E44").Value2 = '  -1.84%  '

# Row 45
Set-TextValue $ws.Range("D45") '1.77'
$ws.Range("E45").Value2 = '  -6.99%  '

# Row 46
$ws.Range("E46").Value2 = '  +0.27%  '

# Row 47
Set-TextValue $ws.Range("D47") '2.38'
$ws.Range("E47").Value2 = '  -7.01%  '

# Row 48
Set-TextValue $ws.Range("D48") '7.09'
$ws.Range("E48").Value2 = '  -1.38%  '

# Row 49
Set-TextValue $ws.Range("D49") '23.15'
$ws.Range("E49").Value2 = '  -3.76%  '

# Row 50
Set-TextValue $ws.Range("D50") '1.11'
$ws.Range("E50").Value2 = '  -5.41%  '

# Row 51
Set-TextValue $ws.Range("D51") '0.914'
$ws.Range("E51").Value2 = '  -4.15%  '
